$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Voice Lines - main")

# Map of row number -> new SnippetID value (column H)
$updates = @{
    2  = "zoSB"
    3  = "zoSB"
    4  = "MrYf"
    5  = "h7gL"
    6  = "jow9"
    7  = "PUQh"
    8  = "LUYG"
    9  = "Gf0M"
    10 = "NiRn"
    11 = "LeCq"
    12 = "QMwJ"
    13 = "sHOs"
    14 = "lGAB"
    15 = "kO4u"
    16 = "gZxI"
    17 = "utPz"
    18 = "utPz"
    19 = "utPz"
    20 = "utPz"
    21 = "utPz"
    22 = "XrkT"
    23 = "4L9m"
    24 = "vdSU"
    25 = "l5L4"
    26 = "iNj4"
    27 = "iNj4"
    28 = "bDdo"
    29 = "0r0f"
}

foreach ($row in $updates.Keys) {
    $ws.Range("H$row").Value = $updates[$row]
}
